# Fix a typo ("achieve an an AUC" -> "achieve an AUC") in the caption
# textbox on slide 7 ("Rectangle 3"). The fix also causes the single run
# that used to hold the whole sentence to be split into three runs, as
# the author re-typed only the "an an " -> "an " portion in the middle
# of the sentence.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(7)
$sh = $s.Shapes.Item("Rectangle 3")
$tr = $sh.TextFrame.TextRange

# --- Step 1: remove the duplicated "an " that causes "an an AUC" ------
$full = $tr.Text
$typo = "achieve an an AUC"
$idx  = $full.IndexOf($typo)
if ($idx -lt 0) {
    throw "Could not locate the typo text to fix"
}

$dupStart0 = $idx + "achieve an ".Length
$dupLen    = "an ".Length
$dupRange  = $tr.Characters($dupStart0 + 1, $dupLen)
if ($dupRange.Text -ne "an ") {
    throw "Unexpected text at duplicate-'an' offset: [$($dupRange.Text)]"
}
$dupRange.Text = ""

# --- Step 2: split the (now-corrected) run into three runs, matching --
# --- the run boundaries the author ended up with -----------------------
$fixed = $tr.Text

$segB = "achieve an "
$segC = "AUC (Area under the Curve) of "

$segBStart0 = $fixed.IndexOf($segB)
$segCStart0 = $fixed.IndexOf($segC)
if ($segBStart0 -lt 0 -or $segCStart0 -lt 0) {
    throw "Could not locate run-split anchors after typo fix"
}

# Re-assigning each segment's own text (to itself) forces the run to be
# split off from its neighbours without altering any character content.
$rangeB = $tr.Characters($segBStart0 + 1, $segB.Length)
$rangeB.Text = $segB

$rangeC = $tr.Characters($segCStart0 + 1, $segC.Length)
$rangeC.Text = $segC

$expected = "Using a 10-fold Cross-Validated resampling with k = 15, the model was able to achieve an AUC (Area under the Curve) of 0.7223 against our test data set."
if ($tr.Text -ne $expected) {
    throw "Post-edit text does not match expected result: [$($tr.Text)]"
}

Write-Host "Slide 7 caption fixed:" $tr.Text
